$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).NumberFormat = "@"
$ws.Cells.Item(2,4).Value = "275.41"
$ws.Cells.Item(2,5).NumberFormat = "@"
$ws.Cells.Item(2,5).Value = "-2.10%"

# Row 3
$ws.Cells.Item(3,4).NumberFormat = "@"
$ws.Cells.Item(3,4).Value = "27.15"
$ws.Cells.Item(3,5).NumberFormat = "@"
$ws.Cells.Item(3,5).Value = "1.00%"

# Row 4
$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = "4.771"
$ws.Cells.Item(4,5).NumberFormat = "@"
$ws.Cells.Item(4,5).Value = "-3.32%"

# Row 5
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = "0.06312"
$ws.Cells.Item(5,5).NumberFormat = "@"
$ws.Cells.Item(5,5).Value = "-1.25%"

# Row 6
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = "6.936"
$ws.Cells.Item(6,5).NumberFormat = "@"
$ws.Cells.Item(6,5).Value = "-0.59%"

# Row 7
$ws.Cells.Item(7,2).Value = "GateToken"
$ws.Cells.Item(7,3).Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = "3.303"
$ws.Cells.Item(7,5).NumberFormat = "@"
$ws.Cells.Item(7,5).Value = "-1.36%"

# Row 8
$ws.Cells.Item(8,2).Value = "FTXToken"
$ws.Cells.Item(8,3).Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = "1.326"
$ws.Cells.Item(8,5).NumberFormat = "@"
$ws.Cells.Item(8,5).Value = "38.91%"

# Row 9
$ws.Cells.Item(9,2).Value = "MXToken"
$ws.Cells.Item(9,3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = "0.8770"
$ws.Cells.Item(9,5).NumberFormat = "@"
$ws.Cells.Item(9,5).Value = "-0.77%"

# Row 10
$ws.Cells.Item(10,2).Value = "WazirX"
$ws.Cells.Item(10,3).Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = "0.1515"
$ws.Cells.Item(10,5).NumberFormat = "@"
$ws.Cells.Item(10,5).Value = "1.28%"

# Row 11
$ws.Cells.Item(11,2).Value = "LiechtensteinCryptoassetsExchange"
$ws.Cells.Item(11,3).Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = "0.05032"
$ws.Cells.Item(11,5).NumberFormat = "@"
$ws.Cells.Item(11,5).Value = "-3.77%"

# Row 12
$ws.Cells.Item(12,2).Value = "MandalaExchangeToken"
$ws.Cells.Item(12,3).Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = "0.07461"
$ws.Cells.Item(12,5).NumberFormat = "@"
$ws.Cells.Item(12,5).Value = "-0.03%"

# Row 13
$ws.Cells.Item(13,2).Value = "BitrueCoin"
$ws.Cells.Item(13,3).Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = "0.02867"
$ws.Cells.Item(13,5).NumberFormat = "@"
$ws.Cells.Item(13,5).Value = "-7.98%"

# Row 14
$ws.Cells.Item(14,2).Value = "BitMartToken"
$ws.Cells.Item(14,3).Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = "0.09024"
$ws.Cells.Item(14,5).NumberFormat = "@"
$ws.Cells.Item(14,5).Value = "-0.04%"

# Row 15
$ws.Cells.Item(15,2).Value = "BitForexToken"
$ws.Cells.Item(15,3).Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = "0.001580"
$ws.Cells.Item(15,5).NumberFormat = "@"
$ws.Cells.Item(15,5).Value = "-0.27%"

# Row 16
$ws.Cells.Item(16,2).Value = "One"
$ws.Cells.Item(16,3).Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = "0.0006369"
$ws.Cells.Item(16,5).NumberFormat = "@"
$ws.Cells.Item(16,5).Value = "0.71%"

# Row 17
$ws.Cells.Item(17,2).Value = "TigerCash"
$ws.Cells.Item(17,3).Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = "0.005839"
$ws.Cells.Item(17,5).NumberFormat = "@"
$ws.Cells.Item(17,5).Value = "-2.48%"

# Row 18
$ws.Cells.Item(18,2).Value = "LEO"
$ws.Cells.Item(18,3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value = "3.449"
$ws.Cells.Item(18,5).NumberFormat = "@"
$ws.Cells.Item(18,5).Value = "-1.63%"

# Row 19
$ws.Cells.Item(19,5).NumberFormat = "@"
$ws.Cells.Item(19,5).Value = "-1.15%"

# Row 20
$ws.Cells.Item(20,5).NumberFormat = "@"
$ws.Cells.Item(20,5).Value = "0.78%"

# Row 21
$ws.Cells.Item(21,5).NumberFormat = "@"
$ws.Cells.Item(21,5).Value = "2.76%"

# Row 22
$ws.Cells.Item(22,5).NumberFormat = "@"
$ws.Cells.Item(22,5).Value = "-0.70%"

# Row 23
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = "0.04409"
$ws.Cells.Item(23,5).NumberFormat = "@"
$ws.Cells.Item(23,5).Value = "1.77%"

# Row 24
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = "0.001172"
$ws.Cells.Item(24,5).NumberFormat = "@"
$ws.Cells.Item(24,5).Value = "-0.14%"

# Row 25
$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = "0.003837"
$ws.Cells.Item(25,5).NumberFormat = "@"
$ws.Cells.Item(25,5).Value = "4.64%"

# Row 26
$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = "0.0001201"
$ws.Cells.Item(26,5).NumberFormat = "@"
$ws.Cells.Item(26,5).Value = "0.43%"

# Row 27
$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = "0.0001939"
$ws.Cells.Item(27,5).NumberFormat = "@"
$ws.Cells.Item(27,5).Value = "14.74%"

# Row 40
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = "0.04094"
$ws.Cells.Item(40,5).NumberFormat = "@"
$ws.Cells.Item(40,5).Value = "-0.19%"

# Row 41
$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = "0.006820"
$ws.Cells.Item(41,5).NumberFormat = "@"
$ws.Cells.Item(41,5).Value = "2.68%"

# Row 42
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = "0.1172"
$ws.Cells.Item(42,5).NumberFormat = "@"
$ws.Cells.Item(42,5).Value = "-0.28%"

# Row 43
$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = "0.002322"
$ws.Cells.Item(43,5).NumberFormat = "@"
$ws.Cells.Item(43,5).Value = "-1.26%"

# Row 44
$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = "0.01154"
$ws.Cells.Item(44,5).NumberFormat = "@"
$ws.Cells.Item(44,5).Value = "-8.10%"

# Row 45
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = "0.00005187"
$ws.Cells.Item(45,5).NumberFormat = "@"
$ws.Cells.Item(45,5).Value = "-1.08%"

# Row 46
$ws.Cells.Item(46,2).Value = "CoinbaseStockToken"
$ws.Cells.Item(46,3).Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = "0.02303"
$ws.Cells.Item(46,5).NumberFormat = "@"
$ws.Cells.Item(46,5).Value = "2.61%"

# Row 47
$ws.Cells.Item(47,2).Value = "BOLO"
$ws.Cells.Item(47,3).Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = "1.490"
$ws.Cells.Item(47,5).NumberFormat = "@"
$ws.Cells.Item(47,5).Value = "-36.77%"
